$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.64
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 6.6
$ws.Range("J2").Value = 3.2
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 1.84
$ws.Range("T2").Value = 1.81
$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.19
$ws.Range("Y2").Value = 19
$ws.Range("Z2").Value = 48
$ws.Range("AB2").Value = 9.199999999999999
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AD2").Value = 23
$ws.Range("AF2").Value = 11.5
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 23
$ws.Range("AJ2").Value = 20
$ws.Range("AL2").Value = 46
$ws.Range("F3").Value = 2.02
$ws.Range("G3").Value = 2.04
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 4
$ws.Range("P3").Value = 2.14
$ws.Range("Q3").Value = 1.71
$ws.Range("R3").Value = 1.46
$ws.Range("S3").Value = 2.78
$ws.Range("U3").Value = 2.26
$ws.Range("V3").Value = 1.35
$ws.Range("W3").Value = 1.96
$ws.Range("AF3").Value = 17.5
$ws.Range("AG3").Value = 11.5
$ws.Range("AK3").Value = 22
$ws.Range("F4").Value = 2.24
$ws.Range("L4").Value = 1.45
$ws.Range("V4").Value = 1.31
$ws.Range("F5").Value = 2.28
$ws.Range("G5").Value = 2.88
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 4.5
$ws.Range("K5").Value = 3.55
$ws.Range("L5").Value = 1.49
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 2.48
$ws.Range("O5").Value = 1.54
$ws.Range("P5").Value = 1.49
$ws.Range("Q5").Value = 2.38
$ws.Range("R5").Value = 1.17
$ws.Range("S5").Value = 4.8
$ws.Range("T5").Value = 1.94
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 1.31
$ws.Range("W5").Value = 1.6
$ws.Range("F6").Value = 3.25
$ws.Range("I6").Value = 2.42
$ws.Range("J6").Value = 2.8
$ws.Range("K6").Value = 3.9
$ws.Range("L6").Value = 1.42
$ws.Range("M6").Value = 1.08
$ws.Range("S6").Value = 3.95
$ws.Range("V6").Value = 1.71
$ws.Range("AC6").Value = 9
